$wb = $excel.ActiveWorkbook

# ===== Sheet 1: Overview =====
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()
$ws1.Cells.Item(2,2).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(2,4).Value = "2016-03-23 17:18:29"
$ws1.Cells.Item(3,2).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(3,4).Value = "2016-03-23 17:16:29"
$ws1.Cells.Item(4,2).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(4,3).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(4,4).Value = "2016-03-23 17:16:29"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/17d81038-38c9-4d0a-a2b5-9b90647a1c6b.md", "", "", "17d81038-38c9-4d0a-a2b5-9b90647a1c6b.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/ffffd6795c51-6fb3-43a7-b7c7-a7226ad88f57.md", "", "", "ffffd6795c51-6fb3-43a7-b7c7-a7226ad88f57.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/ffffff07d0cd49-d89c-46b8-ad29-4e29b6f8bd6b.md", "", "", "ffffff07d0cd49-d89c-46b8-ad29-4e29b6f8bd6b.md")

# ===== Sheet 2: zh-cn =====
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()
$ws2.Cells.Item(2,2).Value = ".md"
$ws2.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws2.Cells.Item(2,5).Value = "2016-03-23 17:18:24"
$ws2.Cells.Item(2,8).Value = "2016-03-23 17:18:56"
$ws2.Cells.Item(2,10).Value = "Include"
$ws2.Cells.Item(3,2).Value = ".md"
$ws2.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$ws2.Cells.Item(3,5).Value = "2016-03-23 17:16:25"
$ws2.Cells.Item(3,8).Value = "2016-03-23 17:17:06"
$ws2.Cells.Item(3,10).Value = "Include"
$ws2.Cells.Item(4,2).Value = ".md"
$ws2.Cells.Item(4,3).Value = "Handed back: in sync with en-US"
$ws2.Cells.Item(4,5).Value = "2016-03-23 17:16:25"
$ws2.Cells.Item(4,8).Value = "2016-03-23 17:17:06"
$ws2.Cells.Item(4,10).Value = "Include"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/17d81038-38c9-4d0a-a2b5-9b90647a1c6b.md", "", "", "17d81038-38c9-4d0a-a2b5-9b90647a1c6b.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/be1f229145fdb172e4f0427ecd8309de7c8844af/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/17d81038-38c9-4d0a-a2b5-9b90647a1c6b.150fb3fc061825919fdf1cc582986b32f55825aa.zh-cn.xlf", "", "", "17d81038-38c9-4d0a-a2b5-9b90647a1c6b.150fb3fc061825919fdf1cc582986b32f55825aa.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/17d81038-38c9-4d0a-a2b5-9b90647a1c6b.md", "", "", "17d81038-38c9-4d0a-a2b5-9b90647a1c6b.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/17d81038-38c9-4d0a-a2b5-9b90647a1c6b.150fb3fc061825919fdf1cc582986b32f55825aa.zh-cn.xlf", "", "", "17d81038-38c9-4d0a-a2b5-9b90647a1c6b.150fb3fc061825919fdf1cc582986b32f55825aa.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/ffffd6795c51-6fb3-43a7-b7c7-a7226ad88f57.md", "", "", "ffffd6795c51-6fb3-43a7-b7c7-a7226ad88f57.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74b8c21fab3bed07889a34a67cbb8fc69884e8f7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.zh-cn.xlf", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/469ab26ec6756d069301a4232ace188b8b124b39/e2e/f5b2ccdf-10b2-4339-94e9-4b45c024f529.md", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/060e062accac4e3034852a59f4940d16af1fd96c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.zh-cn.xlf", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/ffffff07d0cd49-d89c-46b8-ad29-4e29b6f8bd6b.md", "", "", "ffffff07d0cd49-d89c-46b8-ad29-4e29b6f8bd6b.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74b8c21fab3bed07889a34a67cbb8fc69884e8f7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.zh-cn.xlf", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/469ab26ec6756d069301a4232ace188b8b124b39/e2e/f5b2ccdf-10b2-4339-94e9-4b45c024f529.md", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.md")
$ws2.Hyperlinks.Add($ws2.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/060e062accac4e3034852a59f4940d16af1fd96c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.zh-cn.xlf", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.zh-cn.xlf")

# ===== Sheet 3: de-de =====
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()
$ws3.Cells.Item(2,2).Value = ".md"
$ws3.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws3.Cells.Item(2,5).Value = "2016-03-23 17:18:29"
$ws3.Cells.Item(2,8).Value = "2016-03-23 17:19:08"
$ws3.Cells.Item(2,10).Value = "Include"
$ws3.Cells.Item(3,2).Value = ".md"
$ws3.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$ws3.Cells.Item(3,5).Value = "2016-03-23 17:16:29"
$ws3.Cells.Item(3,8).Value = "2016-03-23 17:17:14"
$ws3.Cells.Item(3,10).Value = "Include"
$ws3.Cells.Item(4,2).Value = ".md"
$ws3.Cells.Item(4,3).Value = "Handed back: in sync with en-US"
$ws3.Cells.Item(4,5).Value = "2016-03-23 17:16:29"
$ws3.Cells.Item(4,8).Value = "2016-03-23 17:17:14"
$ws3.Cells.Item(4,10).Value = "Include"
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/17d81038-38c9-4d0a-a2b5-9b90647a1c6b.md", "", "", "17d81038-38c9-4d0a-a2b5-9b90647a1c6b.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/13da6ce42bd86ab61a9285e6a7480f898b0651ec/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/17d81038-38c9-4d0a-a2b5-9b90647a1c6b.150fb3fc061825919fdf1cc582986b32f55825aa.de-de.xlf", "", "", "17d81038-38c9-4d0a-a2b5-9b90647a1c6b.150fb3fc061825919fdf1cc582986b32f55825aa.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/17d81038-38c9-4d0a-a2b5-9b90647a1c6b.md", "", "", "17d81038-38c9-4d0a-a2b5-9b90647a1c6b.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/17d81038-38c9-4d0a-a2b5-9b90647a1c6b.150fb3fc061825919fdf1cc582986b32f55825aa.de-de.xlf", "", "", "17d81038-38c9-4d0a-a2b5-9b90647a1c6b.150fb3fc061825919fdf1cc582986b32f55825aa.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/ffffd6795c51-6fb3-43a7-b7c7-a7226ad88f57.md", "", "", "ffffd6795c51-6fb3-43a7-b7c7-a7226ad88f57.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/58cabb5003088dee0d8b02e530f43e5afe595776/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.de-de.xlf", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8a0db3aa8b50c235d16d26a5444834695dd4ecde/e2e/f5b2ccdf-10b2-4339-94e9-4b45c024f529.md", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ec9663cd88567d2326e968e5806457d05aff63d2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.de-de.xlf", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4c10c1b5c45f11790e807fe8422c0894c66ece8f/e2e/ffffff07d0cd49-d89c-46b8-ad29-4e29b6f8bd6b.md", "", "", "ffffff07d0cd49-d89c-46b8-ad29-4e29b6f8bd6b.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/58cabb5003088dee0d8b02e530f43e5afe595776/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.de-de.xlf", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8a0db3aa8b50c235d16d26a5444834695dd4ecde/e2e/f5b2ccdf-10b2-4339-94e9-4b45c024f529.md", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.md")
$ws3.Hyperlinks.Add($ws3.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ec9663cd88567d2326e968e5806457d05aff63d2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.de-de.xlf", "", "", "f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.de-de.xlf")